$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: headers for second classification table
$ws.Range("B21").Value = "Number of employees"
$ws.Range("B21").Style = "title"
$ws.Range("C21").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("C21").Style = "title"
$ws.Range("D21").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("D21").Style = "title"

# Row 22: Micro
$ws.Range("A22").Value = "Micro"
$ws.Range("A22").Style = "Normal"
$ws.Range("B22").Value = ""
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = ""
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = ""
$ws.Range("D22").Style = "Normal"

# Row 23: Small
$ws.Range("A23").Value = "Small"
$ws.Range("A23").Style = "Normal"
$ws.Range("B23").Value = ""
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = ""
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = ""
$ws.Range("D23").Style = "Normal"

# Row 24: Medium
$ws.Range("A24").Value = "Medium"
$ws.Range("A24").Style = "Normal"
$ws.Range("B24").Value = "=<200 all sectors"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = ""
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "<S$100,000,000"
$ws.Range("D24").Style = "Normal"

# Row 25: Large
$ws.Range("A25").Value = "Large"
$ws.Range("A25").Style = "Normal"
$ws.Range("B25").Value = ">200"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = ""
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "> S$ 100,000,000"
$ws.Range("D25").Style = "Normal"

# Row 32: DSS label (title style)
$ws.Range("A32").Value = "DSS"
$ws.Range("A32").Style = "title"

# Row 33: source text (source/italic style)
$ws.Range("A33").Value = "Department of Statistics Singapore.  Subject: Companies and Businesses, Topic: Enterprises, Title: M600981 - Topline Estimates For All Enterprises And SMEs, Annual.  Utilizing SingStat Table Builder. Singapore."
$ws.Range("A33").Style = "source"

# Remove old rows 26/27 content (previously held DSS + source text)
$ws.Range("A26").Value = $null
$ws.Range("A27").Value = $null
